$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sheet title and header label to reflect the new "through" date
$ws.Name = "Through 2022-07-12"
$ws.Range("B1").Value = "July 2022 (through July 12)"

# Update carjacking counts for neighborhoods/months affected by the 2022-07-20 data refresh
$ws.Range("B2").Value = 7
$ws.Range("W3").Value = 1
$ws.Range("P5").Value = 3
$ws.Range("AR5").Value = 3
$ws.Range("B6").Value = 6
$ws.Range("E7").Value = 4
$ws.Range("AK7").Value = 4
$ws.Range("P8").Value = 8
$ws.Range("AK15").Value = 1
$ws.Range("AR16").Value = 1
$ws.Range("I19").Value = 2
$ws.Range("AY34").Value = 1
$ws.Range("B35").Value = 1
$ws.Range("I38").Value = 2
$ws.Range("AR58").Value = 1
$ws.Range("P65").Value = 1
$ws.Range("AD68").Value = 1
$ws.Range("B70").Value = 1
$ws.Range("AK75").Value = 1
$ws.Range("P94").Value = 1
